$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Population ID" column inserted after Sample (column B); existing
# locus columns shift right by one, and two new loci (locus_05, locus_06)
# are appended at the end.
$ws.Columns.Item(2).Insert()

$ws.Range("B1").Value = "Population ID"
$ws.Range("G1").Value = "locus_05"
$ws.Range("H1").Value = "locus_06"

$populations = @("pop1","pop1","pop1","pop1","pop2","pop2","pop2","pop2","pop3","pop3","pop3","pop3")
for ($i = 0; $i -lt $populations.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $populations[$i]
}

# locus_05 is monomorphic (CC) for every sample.
$ws.Range("G2:G13").Value = "CC"

# locus_06 is monomorphic (TT) except for the last population, which
# segregates a minor allele.
$ws.Range("H2:H10").Value = "TT"
$ws.Range("H11").Value = "TA"
$ws.Range("H12").Value = "TA"
$ws.Range("H13").Value = "TC"

$ws.Columns.Item(1).ColumnWidth = 7.83
$ws.Columns.Item(2).ColumnWidth = 10.83

[void]$ws.Range("H14").Select()
